$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1) Status text for the 5945e1c3... row changed from "Ready for handoff" to
#    "Handback transform failed" everywhere that shared string is used
#    (Overview!E4 / Overview!F4 and zh-cn!C4 / de-de!C4 all reference the same
#    shared string, so every one of these cells must be updated together).
$wsOverview.Range("E4").Value = "Handback transform failed"
$wsOverview.Range("F4").Value = "Handback transform failed"
$wsZhCn.Range("C4").Value = "Handback transform failed"
$wsDeDe.Range("C4").Value = "Handback transform failed"

# 2) Error Detail column (R) widened to 40 characters-equivalent raw width on
#    both the zh-cn and de-de sheets. Excel reports this raw xlsx width of 40
#    as a ColumnWidth of 39.17 (same as the workbook's other width=40 columns).
$wsZhCn.Columns.Item(18).ColumnWidth = 39.17
$wsDeDe.Columns.Item(18).ColumnWidth = 39.17

# 3) Populate the previously-empty Error Detail (R4) cells with the handback
#    mismatch diagnostic message for each locale.
$wsZhCn.Range("R4").Value = "Handback file name: nsm3utas.ydp is different with handoff file name: 5945e1c3-d233-4e7f-9ed4-d4acaf928a8c.c5da13ef5c325faa831b9df6f26b03080b90108d.zh-cn."
$wsDeDe.Range("R4").Value = "Handback file name: nsm3utas.ydp is different with handoff file name: 5945e1c3-d233-4e7f-9ed4-d4acaf928a8c.c5da13ef5c325faa831b9df6f26b03080b90108d.de-de."
